$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.245.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.907.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'307.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.57%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.5251"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.55%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.66%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07307"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.78%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.31%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.54%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08091"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.86%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'96.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.24%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.372"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'1.764.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -7.35%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008695"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.84%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'14.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.68%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.276.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.37%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.128"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.26%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.505"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.348"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'150.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'18.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.30%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'116.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.77%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.857"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.99%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.880"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09226"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.28%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.8211"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.14%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05084"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.235"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'2.758"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.05%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.371"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5753"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.02004"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.45%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.084"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.85%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'9.056"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.613"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'117.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.1524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.52%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4947"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.79%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'10.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.04%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.21%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'38.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.28%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'64.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.27%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05974"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.71%  "
$ws.Range("E51").Style = "Normal"
